$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# The COM layer stores Left/Top (in points) as a 32-bit float, then converts
# back to EMU by truncating (floor). A naive EMU/12700.0 value can therefore
# "round trip" to one EMU less than intended. This helper searches for a
# point value whose float32 round-trip reproduces the exact target EMU.
function EmuToPt($targetEmu) {
    $base = $targetEmu / 12700.0
    $eps = 0.0
    $step = 0.000001
    for ($i = 0; $i -lt 500000; $i++) {
        $candidate = $base + $eps
        $f32 = [float]$candidate
        $val = [int64]([double]$f32 * 12700.0)
        if ($val -eq $targetEmu) {
            return $candidate
        }
        if ($val -gt $targetEmu) {
            break
        }
        $eps += $step
    }
    return $base
}

# Straight connectors that shift right (x: 4541447 -> 4550591), y unchanged
$shiftedConnectors = @(
    "Straight Connector 5",
    "Straight Connector 9",
    "Straight Connector 10",
    "Straight Connector 11",
    "Straight Connector 12",
    "Straight Connector 13",
    "Straight Connector 14",
    "Straight Connector 15",
    "Straight Connector 18",
    "Straight Connector 19",
    "Straight Connector 97",
    "Straight Connector 105"
)
foreach ($name in $shiftedConnectors) {
    $shp = $s.Shapes.Item($name)
    $shp.Left = EmuToPt(4550591)
}

# Rotated connector (rot=5400000): x changes 5065646 -> 5734434, y unchanged
$conn40 = $s.Shapes.Item("Straight Connector 40")
$conn40.Left = EmuToPt(5734434)

# Groups: only the vertical (Top) position changes, x unchanged
$grp65 = $s.Shapes.Item("Group 65")
$grp65.Top = EmuToPt(2390170)

$grp66 = $s.Shapes.Item("Group 66")
$grp66.Top = EmuToPt(3534026)

# TextBox 87 shifts right (4169537 -> 4178681)
$tb87 = $s.Shapes.Item("TextBox 87")
$tb87.Left = EmuToPt(4178681)

# TextBoxes shifting right by 9144 EMU (4188797 -> 4197941)
$shiftedTextboxes9144 = @(
    "TextBox 88",
    "TextBox 89",
    "TextBox 90",
    "TextBox 91",
    "TextBox 92",
    "TextBox 93",
    "TextBox 94",
    "TextBox 103"
)
foreach ($name in $shiftedTextboxes9144) {
    $shp = $s.Shapes.Item($name)
    $shp.Left = EmuToPt(4197941)
}

# TextBoxes shifting right by 9144 EMU (4062159 -> 4071303)
$shiftedTextboxesB = @(
    "TextBox 95",
    "TextBox 96"
)
foreach ($name in $shiftedTextboxesB) {
    $shp = $s.Shapes.Item($name)
    $shp.Left = EmuToPt(4071303)
}

# TextBox 98: big shift right (5018441 -> 5687229)
$tb98 = $s.Shapes.Item("TextBox 98")
$tb98.Left = EmuToPt(5687229)

# Isosceles Triangle 1: big shift right (5229177 -> 5897965)
$tri1 = $s.Shapes.Item("Isosceles Triangle 1")
$tri1.Left = EmuToPt(5897965)
